$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.10%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'0.94%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.703"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.39%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06202"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.81%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'0.55%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8499"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.14%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9161"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.02%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'1.16%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.04651"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-6.15%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07083"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.65%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03146"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.56%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09046"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.94%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001529"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.77%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006169"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.57%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006068"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.53%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.466"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.12%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.169"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.02%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.27%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.44%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'0.89%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.086"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.86%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04250"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.41%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001211"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.56%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-5.82%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.01%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'5.01%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E40").Value = "'1.94%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.24%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.004134"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'4.84%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-0.80%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-8.82%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005174"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.78%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.03591"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-34.21%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.1677"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'23.92%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").Style = "Normal"
